# Applies the k8s-ref-arch-figures.pptx content edit:
#  1. Footer "datetimeFigureOut" field text 10/01/2020 -> 13/05/2021
#     on the slide master and every slide layout.
#  2. Figure 1-1 text relabelling on slide 3 and slide 4
#     (Physical/virtual infra box, infra-mgmt box, Worker Node box,
#      Master->Control Node box, Master->Control OS box,
#      Master->Control Node Services box).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Date placeholder text, master + all custom layouts
# ---------------------------------------------------------------------
$newDate = "13/05/2021"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

$layouts = $master.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    $layout = $layouts.Item($j)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---------------------------------------------------------------------
# 2. Figure text relabelling — identical edits on slide 3 and slide 4
# ---------------------------------------------------------------------
foreach ($slideIdx in 3, 4) {
    $s = $p.Slides.Item($slideIdx)

    # Shape 2: "Physical / virtual compute, storage and network infrastructure"
    #          "(aligns with NFVI)"
    #       -> "Physical / virtual compute, storage"
    #          "and network hardware resources"
    $shp = $s.Shapes.Item(2)
    $tr = $shp.TextFrame.TextRange
    $tr.Paragraphs(1, 1).Text = "Physical / virtual compute, storage"
    $tr.Paragraphs(2, 1).Text = "and network hardware resources"

    # Shape 3: "Virtual or physical infrastructure management (largely aligns with VIM)"
    #       -> "Virtual or physical hardware infrastructure manager"
    $shp = $s.Shapes.Item(3)
    $shp.TextFrame.TextRange.Text = "Virtual or physical hardware infrastructure manager"

    # Shape 4: "Kubernetes Worker Node Machine (Virtual / Physical)"
    #       -> "Kubernetes Worker Node" / "(Virtual / Physical)" (two paragraphs)
    $shp = $s.Shapes.Item(4)
    $shp.TextFrame.TextRange.Text = "Kubernetes Worker Node`r(Virtual / Physical)"

    # Shape 8: "Kubernetes Master Machine (Virtual / Physical)"
    #       -> "Kubernetes Control Node" / "(Virtual / Physical)"
    $shp = $s.Shapes.Item(8)
    $shp.TextFrame.TextRange.Text = "Kubernetes Control Node`r(Virtual / Physical)"

    # Shape 10: "Kubernetes Master OS" -> "Kubernetes Control Node OS" (slide 3)
    #                                  -> "Kubernetes Control OS" (slide 4)
    $shp = $s.Shapes.Item(10)
    if ($slideIdx -eq 3) {
        $shp.TextFrame.TextRange.Text = "Kubernetes Control Node OS"
    } else {
        $shp.TextFrame.TextRange.Text = "Kubernetes Control OS"
    }

    # Shape 11: first run "Kubernetes Master Node Services (" -> "Kubernetes Control Node Services ("
    #           (the remaining runs - "kube-apiserver", ", controller-managers, DNS, CNI, etc.)" - are untouched)
    $shp = $s.Shapes.Item(11)
    $tr = $shp.TextFrame.TextRange
    $firstRunLen = ("Kubernetes Master Node Services (").Length
    $tr.Characters(1, $firstRunLen).Text = "Kubernetes Control Node Services ("
}
